$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.556.16"
$ws.Range("E2").Value = "  -0.69%  "
$ws.Range("D3").Value = "2.447.06"
$ws.Range("E3").Value = "  +0.03%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'508.89"
$ws.Range("E5").Value = "  -3.00%  "
$ws.Range("D6").Value = "'132.90"
$ws.Range("E6").Value = "  +1.47%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("D8").Value = "'0.556"
$ws.Range("E8").Value = "  -1.32%  "
$ws.Range("D9").Value = "2.448.01"
$ws.Range("E9").Value = "  -0.11%  "
$ws.Range("E11").Value = "  -1.06%  "
$ws.Range("E12").Value = "  -0.29%  "
$ws.Range("E13").Value = "  -7.40%  "
$ws.Range("D14").Value = "2.884.36"
$ws.Range("E14").Value = "  +0.05%  "
$ws.Range("D15").Value = "57.541.85"
$ws.Range("E15").Value = "  -0.54%  "
$ws.Range("D16").Value = "'21.87"
$ws.Range("E16").Value = "  +0.62%  "
$ws.Range("E17").Value = "  +0.92%  "
$ws.Range("D18").Value = "2.472.44"
$ws.Range("E18").Value = "  +0.98%  "
$ws.Range("D19").Value = "'10.29"
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("D20").Value = "'4.11"
$ws.Range("E20").Value = "  +0.13%  "
$ws.Range("D21").Value = "'314.11"
$ws.Range("E21").Value = "  +0.97%  "
$ws.Range("E22").Value = "  +5.61%  "
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").Value = "'5.71"
$ws.Range("E24").Value = "  -2.16%  "
$ws.Range("D25").Value = "'65.14"
$ws.Range("E25").Value = "  +0.33%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("D27").Value = "2.549.47"
$ws.Range("E27").Value = "  -1.08%  "
$ws.Range("E28").Value = "  -0.70%  "
$ws.Range("E29").Value = "  -5.27%  "
$ws.Range("D30").Value = "'7.54"
$ws.Range("E30").Value = "  +4.06%  "
$ws.Range("D31").Value = "'173.46"
$ws.Range("E31").Value = "  -0.31%  "
$ws.Range("D32").Value = "0.0₃0732"
$ws.Range("E32").Value = "  -0.31%  "
$ws.Range("E33").Value = "  -0.17%  "
$ws.Range("E34").Value = "  -0.80%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D37").Value = "'0.996"
$ws.Range("E37").Value = "  -0.17%  "
$ws.Range("D38").Value = "'18.00"
$ws.Range("E39").Value = "  +4.56%  "
$ws.Range("D40").Value = "'3.84"
$ws.Range("E40").Value = "  +1.40%  "
$ws.Range("D41").Value = "'36.70"
$ws.Range("E41").Value = "  +1.30%  "
$ws.Range("D42").Value = "'0.813"
$ws.Range("E42").Value = "  -0.48%  "
$ws.Range("E43").Value = "  +0.90%  "
$ws.Range("D44").Value = "'136.13"
$ws.Range("E44").Value = "  +11.87%  "
$ws.Range("E45").Value = "  -0.10%  "
$ws.Range("D46").Value = "'4.89"
$ws.Range("E46").Value = "  +2.07%  "
$ws.Range("D47").Value = "'256.32"
$ws.Range("E47").Value = "  -2.02%  "
$ws.Range("E48").Value = "  -2.28%  "
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("D51").Value = "'0.0214"
$ws.Range("E51").Value = "  +1.17%  "
